$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $c = $ws.Range($cellAddr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "29.174.54"
Set-TextValue "E2" "  -2.84%  "
Set-TextValue "D3" "1.847.76"
Set-TextValue "E3" "  -1.96%  "
Set-TextValue "D4" "1.001"
Set-TextValue "E4" "  +0.05%  "
Set-TextValue "D5" "0.7043"
Set-TextValue "E5" "  -4.49%  "
Set-TextValue "E6" "  -1.32%  "
Set-TextValue "D7" "1.001"
Set-TextValue "E7" "  -0.05%  "
Set-TextValue "D8" "0.3048"
Set-TextValue "E8" "  -3.58%  "
Set-TextValue "D9" "0.07392"
Set-TextValue "E9" "  +3.15%  "
Set-TextValue "D10" "23.43"
Set-TextValue "E10" "  -4.90%  "
Set-TextValue "D11" "0.08135"
Set-TextValue "E11" "  -2.18%  "
Set-TextValue "B12" "WrappedEther"
Set-TextValue "C12" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D12" "1.855.02"
Set-TextValue "E12" "  -3.83%  "
Set-TextValue "B13" "Polygon"
Set-TextValue "C13" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D13" "0.7268"
Set-TextValue "E13" "  -3.84%  "
Set-TextValue "D14" "5.219"
Set-TextValue "E14" "  -3.38%  "
Set-TextValue "D15" "88.80"
Set-TextValue "E15" "  -4.03%  "
Set-TextValue "D16" "29.086.55"
Set-TextValue "E16" "  -3.27%  "
Set-TextValue "D17" "5.758"
Set-TextValue "E17" "  -6.33%  "
Set-TextValue "D18" "238.44"
Set-TextValue "E18" "  -4.44%  "
Set-TextValue "D19" "13.06"
Set-TextValue "E19" "  -3.61%  "
Set-TextValue "D20" "0.000007638"
Set-TextValue "E20" "  -2.68%  "
Set-TextValue "E21" "  +0.00%  "
Set-TextValue "B22" "BinanceUSD"
Set-TextValue "C22" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D22" "1.001"
Set-TextValue "E22" "  +0.10%  "
Set-TextValue "B23" "WrappedliquidstakedEther2.0"
Set-TextValue "C23" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D23" "2.079.87"
Set-TextValue "E23" "  -4.57%  "
Set-TextValue "D24" "7.589"
Set-TextValue "E24" "  -3.93%  "
Set-TextValue "D25" "8.997"
Set-TextValue "E25" "  -2.86%  "
Set-TextValue "D26" "160.59"
Set-TextValue "E26" "  -2.13%  "
Set-TextValue "D27" "0.1452"
Set-TextValue "E27" "  -7.58%  "
Set-TextValue "E28" "  -3.19%  "
Set-TextValue "D29" "1.968"
Set-TextValue "E29" "  -3.76%  "
Set-TextValue "E30" "  -4.99%  "
Set-TextValue "D31" "4.523"
Set-TextValue "E31" "  -0.55%  "
Set-TextValue "D32" "1.492"
Set-TextValue "E32" "  -2.62%  "
Set-TextValue "D33" "3.989"
Set-TextValue "E33" "  -4.63%  "
Set-TextValue "D34" "0.05195"
Set-TextValue "E34" "  -2.29%  "
Set-TextValue "D35" "1.186"
Set-TextValue "E35" "  -4.96%  "
Set-TextValue "D36" "1.029"
Set-TextValue "E36" "  +3.03%  "
Set-TextValue "D37" "0.7038"
Set-TextValue "E37" "  -8.18%  "
Set-TextValue "D38" "2.663"
Set-TextValue "E38" "  -2.43%  "
Set-TextValue "E39" "  -4.47%  "
Set-TextValue "B40" "TrustWalletToken"
Set-TextValue "C40" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D40" "0.9636"
Set-TextValue "E40" "  +9.70%  "
Set-TextValue "B41" "MXToken"
Set-TextValue "C41" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D41" "2.678"
Set-TextValue "E41" "  -2.88%  "
Set-TextValue "D42" "6.010"
Set-TextValue "E42" "  -0.47%  "
Set-TextValue "B43" "TheSandbox"
Set-TextValue "C43" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D43" "0.4292"
Set-TextValue "E43" "  -5.72%  "
Set-TextValue "B44" "Maker"
Set-TextValue "C44" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D44" "1.070.27"
Set-TextValue "E44" "  -1.52%  "
Set-TextValue "D45" "70.36"
Set-TextValue "E45" "  -2.70%  "
Set-TextValue "D46" "1.000"
Set-TextValue "E46" "  -0.09%  "
Set-TextValue "D47" "102.69"
Set-TextValue "D48" "1.740"
Set-TextValue "E48" "  -6.07%  "
Set-TextValue "B49" "Aptos"
Set-TextValue "C49" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D49" "7.047"
Set-TextValue "E49" "  -6.41%  "
Set-TextValue "B50" "RocketPoolETH"
Set-TextValue "C50" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D50" "1.976.55"
Set-TextValue "E50" "  -4.18%  "
Set-TextValue "D51" "9.114"
Set-TextValue "E51" "  -4.31%  "
